$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First fill the whole A2:C7 block with a single space (this registers " " as
# the first newly-added shared string).
$ws.Range("A2:C7").Value = " "

# Then go back and change a few specific cells to a double space, which
# registers "  " as the second newly-added shared string.
$ws.Range("A2").Value = "  "
$ws.Range("A3").Value = "  "
$ws.Range("B3").Value = "  "

$ws.Range("C7").Select()
